$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 389 - everything from row 389 down shifts to row 390+
$ws.Rows(389).Insert()

# Populate the newly inserted row 389 with the new price-record data
$ws.Range("A389").Value = 3
$ws.Range("B389").Value = "Femacal de La Calera"
$ws.Range("C389").Value = "Coquimbo"
$ws.Range("D389").Value = 44943
$ws.Range("E389").Value = 5
$ws.Range("F389").Value = 100112043
$ws.Range("G389").Value = "Pepino ensalada"
$ws.Range("H389").Value = "Sin especificar"
$ws.Range("I389").Value = "Primera"
$ws.Range("J389").Value = 125
$ws.Range("K389").Value = 13500
$ws.Range("L389").Value = 14000
$ws.Range("M389").Value = 13740
$ws.Range("N389").Value = "$/caja 60 unidades"
$ws.Range("O389").Value = "Limache"
$ws.Range("P389").Value = 229
$ws.Range("Q389").Value = 60
$ws.Range("R389").Value = "Hortaliza"
